# Apply test-data tweaks to Sample_Data workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A = "Operation" -- change a few rows from "Add" to "Change" / "Skip"
# ("Skip" is written before "Change" so new shared-string entries land in the
#  same order as the reference workbook.)
$ws.Range("A3").Value = "Skip"
$ws.Range("A2").Value = "Change"
$ws.Range("A4").Value = "Change"
$ws.Range("A5").Value = "Change"

# Column E = "Brief Description" -- tweak wording for min/mult qty test rows
$ws.Range("E4").Value = "Order minimum of 250.  Max is left to default."
$ws.Range("E5").Value = "Order in multiples of 5, to a maximum of 1000."

# New numeric values: Min Qty for row 4 (X), Max Qty for row 5 (Y)
$ws.Range("X4").Value = 250
$ws.Range("Y5").Value = 1000

# Update the frozen-pane top-left cell and the active selection in the view
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 13   # column M
$win.ScrollRow = 2

$ws.Range("AA2").Select()
